$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (shifts heuristic..xi from E..L to F..M)
$ws.Range("E1:E3").EntireColumn.Insert()

# New header for inserted column
$ws.Range("E1").Value = "eta"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New "eta" values
$ws.Range("E2").Value = 0.1
$ws.Range("E3").Value = 0.1

# Updated data values (row 2)
$ws.Range("B2").Value = 1604.768326474094
$ws.Range("D2").Value = 0
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = -2320349377.127979
$ws.Range("J2").Value = 1
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 1115970.9

# Updated data values (row 3)
$ws.Range("B3").Value = 1604.768326474094
$ws.Range("D3").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = -2320349377.127979
$ws.Range("J3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 1115970.9
